$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the tiny floating point correction on the existing A7 timestamp
$ws.Range("A7").Value = 45807.39295208333

# Append the new row of data (row 8)
$ws.Range("A8").Value = 45808.39136096802
$ws.Range("B8").Value = "CREATINA MONOHIDRATO EN POLVO"
$ws.Range("C8").Value = "1Kg"
$ws.Range("D8").Value = "15,41€"

# Match the date/time number format used by the other "fecha" cells
$ws.Range("A8").NumberFormat = $ws.Range("A7").NumberFormat
